# Add new "2020" column (X) to the right of the existing "2019" column (W),
# copying the number formatting/style from the corresponding W cells and
# filling in the 2020 values reported for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the whole W4:W16 column into X4:X16 first so the
# new cells inherit the same number format / font / borders as the year
# columns next to them.
$ws.Range("W4:W16").Copy() | Out-Null
$ws.Range("X4:X16").PasteSpecial(-4122) | Out-Null

# Header for the new column
$ws.Range("X4").Value = 2020

# Data values for 2020, row by row
$ws.Range("X5").Value = 45.3
$ws.Range("X6").Value = 48.2
$ws.Range("X7").Value = 43.6
$ws.Range("X8").Value = 48.8
$ws.Range("X9").Value = 41.5
$ws.Range("X10").Value = 49.7
$ws.Range("X11").Value = 46.7
$ws.Range("X12").Value = 36.5
$ws.Range("X13").Value = 29.6
$ws.Range("X14").Value = 54.7
$ws.Range("X15").Value = 51.6
$ws.Range("X16").Value = 47.2

$ws.Range("X4").Select() | Out-Null
